# Auto-generated edit script applying the cryptos.xlsx price/volume refresh diff
# (commit: "Updated cryptos list on Tue Jul 25 14:55:51 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the literal text into the cell without Excel re-interpreting
    # number-like strings (e.g. "1.000", "0.9344") as actual numbers, and
    # without leaving a lasting Text-format style on the cell afterwards.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.178.81"
Set-TextValue $ws.Range("E2") "  +0.21%  "
Set-TextValue $ws.Range("D3") "1.855.60"
Set-TextValue $ws.Range("E3") "  +0.48%  "
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "0.6968"
Set-TextValue $ws.Range("E5") "  +0.71%  "
Set-TextValue $ws.Range("D6") "237.06"
Set-TextValue $ws.Range("E6") "  -0.30%  "
Set-TextValue $ws.Range("E7") "  -0.07%  "
Set-TextValue $ws.Range("D8") "0.07657"
Set-TextValue $ws.Range("E8") "  +0.82%  "
Set-TextValue $ws.Range("D9") "0.3038"
Set-TextValue $ws.Range("E9") "  -0.18%  "
Set-TextValue $ws.Range("D10") "23.19"
Set-TextValue $ws.Range("E10") "  -0.64%  "
Set-TextValue $ws.Range("D11") "0.08155"
Set-TextValue $ws.Range("E11") "  +0.67%  "
Set-TextValue $ws.Range("D12") "1.838.70"
Set-TextValue $ws.Range("E12") "  -1.19%  "
Set-TextValue $ws.Range("D13") "0.7143"
Set-TextValue $ws.Range("E13") "  -0.94%  "
Set-TextValue $ws.Range("D14") "5.139"
Set-TextValue $ws.Range("E14") "  -0.47%  "
Set-TextValue $ws.Range("D15") "89.21"
Set-TextValue $ws.Range("E15") "  +0.31%  "
Set-TextValue $ws.Range("D16") "29.194.58"
Set-TextValue $ws.Range("E16") "  -0.08%  "
Set-TextValue $ws.Range("D17") "5.732"
Set-TextValue $ws.Range("E17") "  -0.59%  "
Set-TextValue $ws.Range("D18") "13.25"
Set-TextValue $ws.Range("E18") "  +1.40%  "
Set-TextValue $ws.Range("B19") "ShibaInu"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D19") "0.000007696"
Set-TextValue $ws.Range("E19") "  +0.06%  "
Set-TextValue $ws.Range("B20") "BitcoinCash"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "236.75"
Set-TextValue $ws.Range("E20") "  -2.13%  "
Set-TextValue $ws.Range("D21") "0.9996"
Set-TextValue $ws.Range("E21") "  +0.03%  "
Set-TextValue $ws.Range("D22") "2.112.93"
Set-TextValue $ws.Range("E22") "  +0.20%  "
Set-TextValue $ws.Range("E23") "  -0.01%  "
Set-TextValue $ws.Range("D24") "7.423"
Set-TextValue $ws.Range("E24") "  -2.50%  "
Set-TextValue $ws.Range("D25") "0.1474"
Set-TextValue $ws.Range("E25") "  +1.61%  "
Set-TextValue $ws.Range("D26") "161.96"
Set-TextValue $ws.Range("E26") "  +0.58%  "
Set-TextValue $ws.Range("D27") "8.977"
Set-TextValue $ws.Range("E27") "  -0.21%  "
Set-TextValue $ws.Range("D28") "17.99"
Set-TextValue $ws.Range("E28") "  -0.42%  "
Set-TextValue $ws.Range("D29") "2.023"
Set-TextValue $ws.Range("E29") "  +4.72%  "
Set-TextValue $ws.Range("D30") "1.415"
Set-TextValue $ws.Range("E30") "  +1.99%  "
Set-TextValue $ws.Range("D31") "4.419"
Set-TextValue $ws.Range("E31") "  +0.21%  "
Set-TextValue $ws.Range("D32") "1.478"
Set-TextValue $ws.Range("E32") "  -1.02%  "
Set-TextValue $ws.Range("D33") "4.004"
Set-TextValue $ws.Range("E33") "  -1.12%  "
Set-TextValue $ws.Range("D34") "0.05184"
Set-TextValue $ws.Range("E34") "  -0.95%  "
Set-TextValue $ws.Range("D35") "1.159"
Set-TextValue $ws.Range("E35") "  -2.26%  "
Set-TextValue $ws.Range("D36") "0.7093"
Set-TextValue $ws.Range("E36") "  -0.13%  "
Set-TextValue $ws.Range("D37") "0.9981"
Set-TextValue $ws.Range("E37") "  -0.61%  "
Set-TextValue $ws.Range("D38") "2.657"
Set-TextValue $ws.Range("E38") "  -0.11%  "
Set-TextValue $ws.Range("D39") "0.01849"
Set-TextValue $ws.Range("E39") "  -0.42%  "
Set-TextValue $ws.Range("D40") "2.718"
Set-TextValue $ws.Range("E40") "  +1.01%  "
Set-TextValue $ws.Range("D41") "0.9344"
Set-TextValue $ws.Range("E41") "  +1.89%  "
Set-TextValue $ws.Range("D42") "1.142.34"
Set-TextValue $ws.Range("E42") "  +9.76%  "
Set-TextValue $ws.Range("D43") "0.4269"
Set-TextValue $ws.Range("E43") "  -0.39%  "
Set-TextValue $ws.Range("D44") "70.70"
Set-TextValue $ws.Range("E44") "  +1.48%  "
Set-TextValue $ws.Range("D45") "5.860"
Set-TextValue $ws.Range("E45") "  -1.51%  "
Set-TextValue $ws.Range("D47") "103.32"
Set-TextValue $ws.Range("E47") "  +0.88%  "
Set-TextValue $ws.Range("D48") "1.795"
Set-TextValue $ws.Range("E48") "  +2.84%  "
Set-TextValue $ws.Range("D49") "2.010.97"
Set-TextValue $ws.Range("E49") "  -0.01%  "
Set-TextValue $ws.Range("D50") "9.136"
Set-TextValue $ws.Range("D51") "6.945"
Set-TextValue $ws.Range("E51") "  -3.58%  "
